$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the prototype screenshot (InlineShape) that sits alone in the
#    very first paragraph of the document, then collapse the run of empty
#    "spacer" paragraphs (the image paragraph + 7 blank tab-stop paragraphs +
#    one trailing empty paragraph = 9 paragraphs) that used to precede the
#    first table down to a single empty paragraph.
# ---------------------------------------------------------------------------
if ($d.InlineShapes.Count -gt 0) {
    $d.InlineShapes.Item(1).Delete()
}

$firstPara = $d.Paragraphs.Item(1)
$lastSpacerPara = $d.Paragraphs.Item(9)
$spacerRange = $d.Range($firstPara.Range.Start, $lastSpacerPara.Range.Start)
$spacerRange.Delete()

# ---------------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from its old spot (a paragraph right
#    after the final table) to the remaining lead paragraph at the very start
#    of the document.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range(0, 0))

# ---------------------------------------------------------------------------
# 3) Drop the stale <w:lastRenderedPageBreak/> markers in front of the
#    "CT NO." / "TESTE No." table header cells - a plain Find & Replace
#    (replacing the text with itself) rewrites those runs and sheds the
#    leftover page-break marker left over from the old pagination.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("CT NO.", $true, $false, $false, $false, $false, $true, 1, $false, "CT NO.", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("TESTE No.", $true, $false, $false, $false, $false, $true, 1, $false, "TESTE No.", 2) | Out-Null
